$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 191, shifting the existing rows 191..294 down to 192..295.
$ws.Rows("191:191").Insert()

# Populate the newly inserted row 191 with the new weekly record.
$ws.Cells.Item(191, 1).Value = 7
$ws.Cells.Item(191, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(191, 3).Value = "Ñuble"
$ws.Cells.Item(191, 4).Value = 44845
$ws.Cells.Item(191, 5).Value = 16
$ws.Cells.Item(191, 6).Value = 100112009
$ws.Cells.Item(191, 7).Value = "Acelga"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 300
$ws.Cells.Item(191, 11).Value = 600
$ws.Cells.Item(191, 12).Value = 700
$ws.Cells.Item(191, 13).Value = 650
$ws.Cells.Item(191, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(191, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(191, 16).Value = 650
$ws.Cells.Item(191, 17).Value = 1
$ws.Cells.Item(191, 18).Value = "Hortaliza"
